# Rename the worksheet "Checklist" -> "Session"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")
$ws.Name = "Session"

# Update the "Type" column (E) values from "Selection" to "Scan" for rows 2-6
for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E = 5
    if ($cell.Value2 -eq "Selection") {
        $cell.Value2 = "Scan"
    }
}
